# Auto-committed on 2022/02/22 週二
# Insert a new "AdminFg" (管理者權限記號) row into the DBD field-list table,
# right after the "Station" (站別) row (old row 39 -> new row 40), pushing
# every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# 1) Insert a new blank row at 40 (shifts rows 40..46 down to 41..47).
$ws.Rows.Item(40).Insert()

# 2) Give the new row the same look & feel as the row above it (the
#    "Station" data row) by copying its formatting into the new row.
$ws.Range("A39:G39").Copy()
$ws.Range("A40:G40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Re-establish the auto-incrementing SEQ formula (column A) for every
#    row from the newly inserted one through the end of the table so the
#    numbering stays continuous (=<row above>+1).
for ($r = 40; $r -le 46; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+1"
}

# 4) Fill in the new row's field-definition data.
$ws.Range("B40").Value = "AdminFg"
$ws.Range("C40").Value = "管理者權限記號"
$ws.Range("D40").Value = "decimal"
$ws.Range("E40").Value = "1"
$ws.Range("G40").Value = "0.否 1.是"

# 5) Restore a sensible selection/view state on the sheet.
$ws.Activate()
$ws.Range("G44").Select()
